$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update populations for Buenos Aires and Toronto
$ws.Range("C4").Value = 2891000.0
$ws.Range("C12").Value = 2800000.0

# Sort the data range (A1:D13) by Population (column C) descending, with headers
$rng = $ws.Range("A1:D13")
$rng.Sort($ws.Range("C1"), 2, $null, $null, 1, $null, 1, 1)

# Apply AutoFilter over the table range
$ws.Range("A1:D13").AutoFilter() | Out-Null

# Register the (normally auto-created) hidden _FilterDatabase defined name
$fdb = $ws.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$A`$1:`$D`$13", $false)
$fdb.Visible = $false
